$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data value changes (date format d/m/yyyy -> yyyymmdd, email update) ---
$ws.Range("C2").Value = "19940315"
$ws.Range("I2").Value = "20170601"
$ws.Range("J2").Value = "20170601"
$ws.Range("M2").Value = "20151220"
$ws.Range("R2").Value = "trannhatphuong@vietbank.com.vn"
